$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Registration sheet: update values, drop the mailto hyperlinks that used
#    to live on A2/A3, and touch the column widths so a second (empty) column
#    definition shows up, matching the new layout.
# ---------------------------------------------------------------------------
$wsReg = $wb.Worksheets.Item("Registration")
$wsReg.Range("A2").Hyperlinks.Delete()
$wsReg.Range("A3").Hyperlinks.Delete()
$wsReg.Range("A2").Value = "ddd"
$wsReg.Range("A3").Value = "lld"
$wsReg.Columns.Item(1).ColumnWidth = 12.1640625
$wsReg.Columns.Item(2).ColumnWidth = 12.83203125

# ---------------------------------------------------------------------------
# 2. Parameter sheet: append a new row of sample data with its own mailto
#    hyperlink, re-using the existing hyperlink cell style.
# ---------------------------------------------------------------------------
$wsParam = $wb.Worksheets.Item("Parameter")
$wsParam.Range("A5").Value = "fsd@kk.com"
$wsParam.Range("B5").Value = "kdkffsfk"
$paramLink = $wsParam.Hyperlinks.Add($wsParam.Range("A5"), "mailto:fsd@kk.com")
$wsParam.Range("A2").Copy()
$wsParam.Range("A5").PasteSpecial(-4122)
$wsParam.Range("A5").Value = "fsd@kk.com"

# ---------------------------------------------------------------------------
# 3. Insert a brand-new "random" worksheet right before test_suite (it will
#    inherit test_suite's old relationship id and become the active tab,
#    same as test_suite used to be).
# ---------------------------------------------------------------------------
$wsSuiteBefore = $wb.Worksheets.Item("test_suite")
$wsRandom = $wb.Worksheets.Add($wsSuiteBefore)
$wsRandom.Name = "random"
$wsRandom.Range("A1").Value = "userName"
$wsRandom.Range("A2").Value = "jfqn@test.com"
$wsRandom.Columns.Item(1).ColumnWidth = 9.3359375

# ---------------------------------------------------------------------------
# 4. test_suite sheet: append a new Registration/N row, copying the font
#    formatting already used by the rows above it. Re-fetch the sheet by
#    name since the handle obtained before the insert now refers to the
#    freshly created "random" sheet instead.
# ---------------------------------------------------------------------------
$wsSuite = $wb.Worksheets.Item("test_suite")
$wsSuite.Range("A5").Copy()
$wsSuite.Range("A6").PasteSpecial(-4122)
$wsSuite.Range("A6").Value = "Registration"
$wsSuite.Range("B6").Value = "N"

$wsRandom.Select()
